# Applies the BoM update for rcbus-opl3:
#  - C2 added to the 10uF polarized-cap group, value "10uF" -> "10uf", qty 3 -> 4
#  - CON1 audio jack "Value" field corrected to "3.5mm jack"
#  - Y1 crystal "Value" field corrected to "14.3181MHz"
#  - Date bumped 2025-10-01 -> 2025-10-04
#  - Component/Fitted/Total counts bumped from 31(28 SMD/3 THT) to 32(29 SMD/3 THT)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# --- Row 12 (BoM row 4): Polarized capacitor group C8 C15 C16 -> C2 C8 C15 C16, 10uF -> 10uf, qty 3 -> 4
# (Quantity Per PCB is stored as text in this sheet, e.g. "3" not 3, so we
# assign a text "4" to keep the same text-typed cell rather than a number.)
$ws.Range("D12").Value = "C2 C8 C15 C16"
$ws.Range("E12").Value = "10uf"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "4"

# --- Row 14 (BoM row 6): Audio jack CON1, Value SJ1-3523N -> 3.5mm jack
$ws.Range("E14").Value = "3.5mm jack"

# --- Row 24 (BoM row 16): Crystal Y1, Value SG-8002CA -> 14.3181MHz
$ws.Range("E24").Value = "14.3181MHz"

# --- Header info block
# Force text format first so Excel doesn't auto-convert the ISO-looking
# date string into a date serial number (the source cell is plain text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2025-10-04"
$ws.Range("F3").Value = "32 (29 SMD/ 3 THT)"
$ws.Range("F4").Value = "32 (29 SMD/ 3 THT)"
$ws.Range("F6").Value = 32
